# Auto-generated PowerShell COM-interop script
# Applies the "Updated symbol list" commit: refreshed Price (D) values
# and two Volume(1h) (E) label corrections in the crypto price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates -------------------------------------------
# These cells hold numeric-looking values stored as TEXT in the workbook
# (e.g. "0.00000000749" must stay literal text, not become 7.49E-09), so we
# force the Text number format before assigning the new value, exactly as
# Excel requires to keep a numeric-looking entry as a text string.
$priceUpdates = @{
    "D2" = "266.38"
    "D3" = "21.39"
    "D4" = "6.266"
    "D5" = "0.06165"
    "D6" = "3.569"
    "D7" = "6.568"
    "D8" = "1.364"
    "D9" = "0.8256"
    "D10" = "0.01341"
    "D11" = "0.1553"
    "D12" = "0.08194"
    "D14" = "0.03210"
    "D15" = "0.09300"
    "D16" = "3.630"
    "D17" = "0.001627"
    "D18" = "0.04692"
    "D19" = "0.006403"
    "D20" = "0.005807"
    "D21" = "0.001067"
    "D22" = "0.0001498"
    "D23" = "3.723"
    "D24" = "2.337"
    "D25" = "0.3294"
    "D26" = "0.1242"
    "D28" = "0.0002723"
    "D40" = "0.04654"
    "D41" = "0.007018"
    "D42" = "0.1141"
    "D43" = "0.003652"
    "D44" = "0.01148"
    "D45" = "0.00005854"
    "D46" = "0.0009886"
    "D47" = "0.00000000749"
    "D48" = "0.7811"
    "D49" = "0.002416"
    "D50" = "0.00002397"
    "D51" = "0.01238"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
}

# --- Column E ("Volume(1h)") label corrections -----------------------------
$volumeUpdates = @{
    "E28" = "27UpBotsUBXTWorstin24h"
    "E49" = "48BOLOBOLO"
}

foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}

